$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ecological_params")
$ws.Range("B1").Value = "PLA_virgin"
$ws.Range("C1").Value = "PLA_recycled"
$ws.Range("D1").Value = "PLA_recycled_industrial"
$ws.Activate()
$ws.Range("B1:D1").Select()
